$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "asas"
$ws.Range("B4").Value = "hafizurrahman.naoe@gmail.com"
$ws.Range("C4").Value = "asasa"

$ws.Range("A5").Value = "sdsd"
$ws.Range("B5").Value = "hafizurrahman.naoe@gmail.com"
$ws.Range("C5").Value = "sdsd"
